# "Added the tasks for today"
#
# Adds a new "Namespace" reference table to the "License Data Ver 2.0" sheet
# (rows 69-72) and re-points the "namespace" column of the existing
# License_Data table (row 7) at it with a new data type + FK constraint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("License Data Ver 2.0")

# --- 1. Re-fix the merged section-header rows first -------------------------
# (B/C/D of every section title row should share the title cell's centered
# style, not the old plain one) - do this before using row 23 as a style
# template for the new section below.
foreach ($r in 23, 29, 35, 41, 54, 62) {
    $ws.Range("A$r").Copy()
    $ws.Range("B$r`:D$r").PasteSpecial(-4122)
}

# --- 2. Build the new "Namespace" table block (rows 69-72) -----------------
# Row 69: merged section title "Namespace" (style copied from another
# section-header row, now fixed above).
$ws.Range("A23:D23").Copy()
$ws.Range("A69:D69").PasteSpecial(-4122)
$ws.Range("A69").Value = "Namespace"
$ws.Range("A69:D69").Merge()

# Row 70: column headers (Column Name / Column Data Type / Constraints / Description)
$ws.Range("A13:D13").Copy()
$ws.Range("A70:D70").PasteSpecial(-4122)
$ws.Range("A70").Value = $ws.Range("A13").Value()
$ws.Range("B70").Value = $ws.Range("B13").Value()
$ws.Range("C70").Value = $ws.Range("C13").Value()
$ws.Range("D70").Value = $ws.Range("D13").Value()

# Row 71: the standard "Id" column row
$ws.Range("A3:D3").Copy()
$ws.Range("A71:D71").PasteSpecial(-4122)
$ws.Range("A71").Value = $ws.Range("A3").Value()
$ws.Range("B71").Value = $ws.Range("B3").Value()
$ws.Range("C71").Value = $ws.Range("C3").Value()
$ws.Range("D71").Value = $ws.Range("D3").Value()

# Row 72: the new "namespaceText" column definition
$ws.Range("A3:D3").Copy()
$ws.Range("A72:D72").PasteSpecial(-4122)
$ws.Range("A72").Value = "namespaceText"
$ws.Range("B72").Value = "Varchar(30)"
$ws.Range("C72").Value = $ws.Range("C26").Value()
$ws.Range("D72").Value = "Namespace Text"

# Row heights / custom-height flag for the new rows, matching the rest of the sheet
$ws.Rows.Item(69).RowHeight = 15.75
$ws.Rows.Item(70).RowHeight = 15.75
$ws.Rows.Item(71).RowHeight = 15.75
$ws.Rows.Item(72).RowHeight = 15.75

# --- 3. Point the "namespace" row of License_Data at the new table ---------
$ws.Range("B7").Value = "Varchar(30)"
$ws.Range("C7").Value = "Foreign key to namespace table"

# --- 4. Misc cosmetic bits from the same editing session -------------------
$ws.Columns.Item(3).ColumnWidth = 27.7
$ws.Activate()
$ws.Range("G7").Select()
